# Atualização de bases das ligas, do dia: 26-04-2024 às 22:13
#
# The underlying data update swaps the full match record (every column
# except the running row index in column A, and the unchanged Div/Date
# columns C and D) between specific pairs of rows on the single
# worksheet "Bolivia Primera División".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose data (columns B, and E through AB) must be swapped
# with each other. Columns A (index), C (Div) and D (Date) are identical
# between each pair and stay untouched.
$rowPairs = @(
    @(27, 28),
    @(47, 48),
    @(104, 105),
    @(107, 108),
    @(128, 129),
    @(142, 145)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Capture the "id" column (B) for both rows.
    $b1 = $ws.Range("B$r1").Value2
    $b2 = $ws.Range("B$r2").Value2

    # Capture the remaining match data (E through AB) for both rows.
    $rest1 = $ws.Range("E$r1`:AB$r1").Value2
    $rest2 = $ws.Range("E$r2`:AB$r2").Value2

    # Swap the captured values between the two rows.
    $ws.Range("B$r1").Value2 = $b2
    $ws.Range("B$r2").Value2 = $b1

    $ws.Range("E$r1`:AB$r1").Value2 = $rest2
    $ws.Range("E$r2`:AB$r2").Value2 = $rest1
}
